$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - ŞANS OYUNLARI
$ws.Range("D2").Value = "33,33 TL - 33,33 TL"
$ws.Range("E2").Value = "26 TL - 26 TL"
$ws.Range("H2").Value = ""

# Row 3 - HESAPTAN EFT - Şube
$ws.Range("C3").Value = "30,46 TL - 60,94 TL - 609,43 TL"

# Row 4 - HESAPTAN EFT - ATM
$ws.Range("C4").Value = "30,46 TL - 60,94 TL - 609,43 TL"

# Row 5 - HESAPTAN EFT - Mobil
$ws.Range("C5").Value = "30,46 TL - 60,94 TL - 609,43 TL"

# Row 6 - DÜZENLİ EFT
$ws.Range("C6").Value = "6,09 TL - 12,19 TL - 152,35 TL"
$ws.Range("D6").Value = ""
$ws.Range("E6").Value = "8.300,01 TL - 199,41 TL"

# Row 7 - KREDİ KARTINDAN FATURA ÖDEME
$ws.Range("H7").Value = ""

# Row 8 - HESAPTAN HAVALE - Şube
$ws.Range("C8").Value = "14,29 TL - 28,57 TL - 300 TL"

# Row 9 - HESAPTAN HAVALE - ATM
$ws.Range("C9").Value = "14,29 TL - 28,57 TL - 300 TL"

# Row 10 - HESAPTAN HAVALE - Mobil
$ws.Range("C10").Value = "14,29 TL - 28,57 TL - 300 TL"

# Row 11 - DÜZENLİ HAVALE
$ws.Range("C11").Value = "3,04 TL - 6,09 TL - 76,17 TL"

# Row 12 - GİDEN SWIFT
$ws.Range("C12").Value = "WU: 1.000,01 USD–9,51 USD"
$ws.Range("D12").Value = ""

# Row 13 - GELEN SWIFT
$ws.Range("C13").Value = "Hesaba: Asgari 0 TL | Azami 9.999.999.999.999 TL"
$ws.Range("D13").Value = ""
$ws.Range("E13").Value = "Hesaba: Asgari 1 TL | Azami 1.114 TL"

# Row 14 - GİDEN SWIFT - Mobil
$ws.Range("C14").Value = "40.000 TL - 1.904,76 TL"
$ws.Range("D14").Value = ""
$ws.Range("E14").Value = "2.170 TL - 2.170 TL"

# Row 15 - ÇEK TAHSİLİ BAŞKA BANKA
$ws.Range("E15").Value = "%0,6 Asgari Tutar: 350 TL Azami Tutar: 350 TL / 15.000 TL"

# Row 17 - AYNI ŞUBE ÇEK TAHSİLATI
$ws.Range("E17").Value = "%0,6 Asgari Tutar: 390 TL Azami Tutar: 390 TL / 4.200 TL"

# Row 20 - ÇEK İADE
$ws.Range("E20").Value = "390 TL"

# Row 21 - BLOKE ÇEK DÜZENLEME
$ws.Range("E21").Value = "%0,5 Asgari Tutar: 1.630 TL Azami Tutar: 1.630 TL / 16.360 TL"

# Row 22 - YP ÇEK TAKASA GÖNDERME
$ws.Range("E22").Value = "%1 Asgari Tutar: 1.170 TL Azami Tutar: 1.170 TL / 29.450 TL"

# Row 23 - ÇEK KARNESİ SAYFA ÜCRETİ
$ws.Range("E23").Value = "75 TL"

# Row 24 - SENET TAHSİLE ALMA
$ws.Range("E24").Value = "780 TL"

# Row 25 - MUAMELESİZ SENET İADESİ
$ws.Range("E25").Value = "780 TL"
